$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($row, $vals)
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

Set-RowValues 2 @("ECs", "Fgf2", "Fgfr4", "ECs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 1, 0.3333333333333333, 0.5555316666666666, 1.666595, 0.01938483203642842, 0.01938483203642843, 0.3590064139122222, 3.23105772521, 0.0006710739118137242, 0.0006710739118137242)
Set-RowValues 3 @("ECs", "Fgf2", "Fgfr4", "FAPs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 2, 0.6666666666666666, 0.1817723333333333, 0.545317, 0.006342799811357313, 0.006342799811357313, 0.1174684315117778, 1.057215883606, 0.0002195782492858341, 0.0002195782492858341)
Set-RowValues 4 @("ECs", "Fgf2", "Fgfr4", "sCs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 27.92075533333333, 83.762266, 0.9742723681522142, 0.9742723681522143, 18.04349031277644, 162.391412814988, 0.03372785320188871, 0.03372785320188871)
Set-RowValues 5 @("FAPs", "Fgf2", "Fgfr4", "ECs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 1, 0.3333333333333333, 0.5555316666666666, 1.666595, 0.01938483203642842, 0.01938483203642843, 8.512315065823334, 76.61083559241001, 0.01591167274022479, 0.01591167274022479)
Set-RowValues 6 @("FAPs", "Fgf2", "Fgfr4", "FAPs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 2, 0.6666666666666666, 0.1817723333333333, 0.545317, 0.006342799811357313, 0.006342799811357313, 2.785265835280667, 25.06739251752601, 0.005206367260000878, 0.005206367260000878)
Set-RowValues 7 @("FAPs", "Fgf2", "Fgfr4", "sCs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 27.92075533333333, 83.762266, 0.9742723681522142, 0.9742723681522143, 427.8248757612387, 3850.423881851148, 0.7997130464039901, 0.7997130464039901)
Set-RowValues 8 @("sCs", "Fgf2", "Fgfr4", "ECs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 1, 0.3333333333333333, 0.5555316666666666, 1.666595, 0.01938483203642842, 0.01938483203642843, 1.499039982953333, 13.49135984658, 0.002802085384389909, 0.002802085384389909)
Set-RowValues 9 @("sCs", "Fgf2", "Fgfr4", "FAPs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 2, 0.6666666666666666, 0.1817723333333333, 0.545317, 0.006342799811357313, 0.006342799811357313, 0.4904922829986667, 4.414430546988001, 0.0009168543020706002, 0.0009168543020706003)
Set-RowValues 10 @("sCs", "Fgf2", "Fgfr4", "sCs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 27.92075533333333, 83.762266, 0.9742723681522142, 0.9742723681522143, 75.34103114240267, 678.0692802816241, 0.1408314685463354, 0.1408314685463354)
